$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 28
$ws.Cells.Item(22, 3).Value = 4
$ws.Cells.Item(22, 4).Value = 28
$ws.Cells.Item(22, 5).Value = 30
$ws.Cells.Item(22, 6).Value = 60
$ws.Cells.Item(22, 7).Value = 90

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 27
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 29
$ws.Cells.Item(23, 5).Value = 34
$ws.Cells.Item(23, 6).Value = 56
$ws.Cells.Item(23, 7).Value = 90
